$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9596.85
$ws.Range("I61").Value = 7405
$ws.Range("K61").Value = 7405
$ws.Range("M61").Value = -7193

$ws.Range("H97").Value = 650.1818
$ws.Range("I97").Value = 437.05264
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 437.05264
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = 58.94736
$ws.Range("N97").Value = -2992

$ws.Range("H122").Value = 6946237.5
$ws.Range("I122").Value = 1702.6666
$ws.Range("J122").Value = 13890772
$ws.Range("K122").Value = 5107.9998
$ws.Range("L122").Value = 41672316
$ws.Range("M122").Value = -2657.9998
$ws.Range("N122").Value = -41677216

$ws.Range("H136").Value = 9596.85
$ws.Range("I136").Value = 7405
$ws.Range("K136").Value = 22215
$ws.Range("M136").Value = -19665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8519.75
$ws.Range("I31").Value = 9921.77
$ws.Range("J31").Value = 5916
$ws.Range("K31").Value = 9921.77
$ws.Range("L31").Value = 5916
$ws.Range("M31").Value = -9626.77
$ws.Range("N31").Value = -6506

$ws.Range("H34").Value = 8519.75
$ws.Range("I34").Value = 9921.77
$ws.Range("J34").Value = 5916
$ws.Range("K34").Value = 9921.77
$ws.Range("L34").Value = 5916
$ws.Range("M34").Value = -9719.77
$ws.Range("N34").Value = -6320

$ws.Range("H58").Value = 1936729.8
$ws.Range("I58").Value = 2934231
$ws.Range("J58").Value = 4070.75
$ws.Range("K58").Value = 2934231
$ws.Range("L58").Value = 4070.75
$ws.Range("M58").Value = -2934028
$ws.Range("N58").Value = -4476.75

$ws.Range("H99").Value = 1884.1818
$ws.Range("I99").Value = 1539.5294
$ws.Range("J99").Value = 3056
$ws.Range("K99").Value = 1539.5294
$ws.Range("L99").Value = 3056
$ws.Range("M99").Value = -41.5293999999999
$ws.Range("N99").Value = -6052

$ws.Range("H126").Value = 1884.1818
$ws.Range("I126").Value = 1539.5294
$ws.Range("J126").Value = 3056
$ws.Range("K126").Value = 4618.5882
$ws.Range("L126").Value = 9168
$ws.Range("M126").Value = -2148.5882
$ws.Range("N126").Value = -14108

$ws.Range("H136").Value = 1936729.8
$ws.Range("I136").Value = 2934231
$ws.Range("J136").Value = 4070.75
$ws.Range("K136").Value = 8802693
$ws.Range("L136").Value = 12212.25
$ws.Range("M136").Value = -8800143
$ws.Range("N136").Value = -17312.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3625944
$ws.Range("I5").Value = 373.7857
$ws.Range("K5").Value = 1121.3571
$ws.Range("M5").Value = -1009.3571

$ws.Range("H12").Value = 47619336
$ws.Range("I12").Value = 111111336
$ws.Range("J12").Value = 338.75
$ws.Range("K12").Value = 333334008
$ws.Range("L12").Value = 1016.25
$ws.Range("M12").Value = -333333835
$ws.Range("N12").Value = -1362.25

$ws.Range("H18").Value = 9091340
$ws.Range("I18").Value = 10000374
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 30001122
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -30000953
$ws.Range("N18").Value = -3338

$ws.Range("H46").Value = 2337.6
$ws.Range("J46").Value = 2337.6
$ws.Range("L46").Value = 7012.799999999999
$ws.Range("N46").Value = -7194.799999999999

$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1034.6666
$ws.Range("L92").Value = 3103.9998
$ws.Range("N92").Value = -5599.9998

$ws.Range("H101").Value = 7174.2
$ws.Range("J101").Value = 7174.2
$ws.Range("L101").Value = 21522.6
$ws.Range("N101").Value = -26390.6

$ws.Range("H103").Value = 1387
$ws.Range("I103").Value = 208
$ws.Range("J103").Value = 3352
$ws.Range("K103").Value = 624
$ws.Range("L103").Value = 10056
$ws.Range("M103").Value = 255
$ws.Range("N103").Value = -11814

$ws.Range("H104").Value = 2137.4
$ws.Range("J104").Value = 2137.4
$ws.Range("L104").Value = 6412.200000000001
$ws.Range("N104").Value = -11654.2

$ws.Range("H106").Value = 3626.3333
$ws.Range("J106").Value = 3626.3333
$ws.Range("L106").Value = 10878.9999
$ws.Range("N106").Value = -12770.9999

$ws.Range("H107").Value = 1067.1538
$ws.Range("I107").Value = 352.25
$ws.Range("J107").Value = 1384.8889
$ws.Range("K107").Value = 1056.75
$ws.Range("L107").Value = 4154.6667
$ws.Range("M107").Value = 863.25
$ws.Range("N107").Value = -7994.6667

$ws.Range("H108").Value = 4725
$ws.Range("I108").Value = 1300
$ws.Range("J108").Value = 15000
$ws.Range("K108").Value = 3900
$ws.Range("L108").Value = 45000
$ws.Range("M108").Value = -1020
$ws.Range("N108").Value = -50760

$ws.Range("H109").Value = 4115
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 4115
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 12345
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -14425

$ws.Range("H110").Value = 4014.0527
$ws.Range("I110").Value = 1207
$ws.Range("J110").Value = 4170
$ws.Range("K110").Value = 3621
$ws.Range("L110").Value = 12510
$ws.Range("M110").Value = 469
$ws.Range("N110").Value = -20690

$ws.Range("H111").Value = 3678.6667
$ws.Range("I111").Value = 1175.6
$ws.Range("J111").Value = 6807.5
$ws.Range("K111").Value = 3526.8
$ws.Range("L111").Value = 20422.5
$ws.Range("M111").Value = -459.7999999999997
$ws.Range("N111").Value = -26556.5

$ws.Range("H112").Value = 2916.3333
$ws.Range("I112").Value = 1419.4
$ws.Range("J112").Value = 3985.5715
$ws.Range("K112").Value = 4258.200000000001
$ws.Range("L112").Value = 11956.7145
$ws.Range("M112").Value = -3150.200000000001
$ws.Range("N112").Value = -14172.7145

$ws.Range("H113").Value = 672.9367999999999
$ws.Range("I113").Value = 682.8
$ws.Range("J113").Value = 635.95
$ws.Range("K113").Value = 2048.4
$ws.Range("L113").Value = 1907.85
$ws.Range("M113").Value = 121.6000000000004
$ws.Range("N113").Value = -6247.85

$ws.Range("H114").Value = 807.46155
$ws.Range("I114").Value = 276.57144
$ws.Range("J114").Value = 1426.8334
$ws.Range("K114").Value = 829.71432
$ws.Range("L114").Value = 4280.5002
$ws.Range("M114").Value = 2424.28568
$ws.Range("N114").Value = -10788.5002

$ws.Range("H135").Value = 3625944
$ws.Range("I135").Value = 373.7857
$ws.Range("K135").Value = 3364.0713
$ws.Range("M135").Value = -829.0713000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5113.5
$ws.Range("I122").Value = 5212.2856
$ws.Range("J122").Value = 4422
$ws.Range("K122").Value = 15636.8568
$ws.Range("L122").Value = 13266
$ws.Range("M122").Value = -13186.8568
$ws.Range("N122").Value = -18166

$ws.Range("H132").Value = 5208.965
$ws.Range("I132").Value = 3827.5227
$ws.Range("J132").Value = 9884.615
$ws.Range("K132").Value = 11482.5681
$ws.Range("L132").Value = 29653.845
$ws.Range("M132").Value = -8952.5681
$ws.Range("N132").Value = -34713.845
